$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("U7:U82").ClearContents()
